# Inventario_SevenShoes1.xlsx — "Cantidad" column (H) tweaks:
#  - H4 and H9 become the text "1" (instead of numeric 2)
#  - the whole quantity column (H2:H10) is re-formatted as a centred
#    integer ("0") instead of the previous text-ish numFmt
#  - selection moves from H12 to I11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Write the new (text) values first, while the column still has its
#    original number format, so Excel stores them as genuine text ("1"),
#    not as the number 1 — matches the shared-string result in the diff.
$ws.Range("H4").Value = "1"
$ws.Range("H9").Value = "1"

# 2) Re-format the quantity column: integer number format, centred.
$qty = $ws.Range("H2:H10")
$qty.HorizontalAlignment = -4108   # xlCenter
$qty.NumberFormat = "0"

# 3) Move the active selection to I11 (was H12).
$ws.Range("I11").Select()
